$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Saudi Professional League: Al-Shabab (KSA) vs Al Najma Club) ---
$ws.Range("I2").Value = 9.800000000000001
$ws.Range("J2").Value = 4.1

# --- Row 3 (Saudi Professional League: Al-Fateh (KSA) vs Al-Kholood Club) ---
$ws.Range("G3").Value = 2.38
$ws.Range("H3").Value = 3.3
$ws.Range("J3").Value = 3.2

# --- Row 4 (Saudi Professional League: Al Ahli vs Al-Khaleej Saihat) ---
$ws.Range("K4").Value = 5.2

# --- Row 5 (English National League: Southend vs Eastleigh) ---
$ws.Range("G5").Value = 1.47
$ws.Range("J5").Value = 4.7

# --- Row 6 (English National League: Yeovil vs Aldershot) ---
$ws.Range("F6").Value = 2.64
$ws.Range("G6").Value = 3
$ws.Range("I6").Value = 3.05

# --- Row 7 (English National League: Boreham Wood vs Scunthorpe) ---
$ws.Range("J7").Value = 3.45

# --- Insert a brand-new row at position 9, pushing the old row 9
#     (Scottish Premiership: Livingston vs St Mirren) down to row 10 ---
$ws.Rows.Item(9).Insert()

# New row 9 holds the former row-8 data (English National League:
# Hartlepool vs Gateshead), except Q9 differs slightly (1.51 vs 1.49).
$ws.Range("A9").Value = "English National League"
$ws.Range("B9").Value = "2026-01-20"
$ws.Range("C9").Value = "16:45:00"
$ws.Range("D9").Value = "Hartlepool"
$ws.Range("E9").Value = "Gateshead"
$ws.Range("F9").Value = 1.56
$ws.Range("G9").Value = 1.71
$ws.Range("H9").Value = 4.7
$ws.Range("I9").Value = 6.4
$ws.Range("J9").Value = 4.3
$ws.Range("K9").Value = 5.1
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 2.4
$ws.Range("Q9").Value = 1.51
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 0
$ws.Range("Y9").Value = 0
$ws.Range("Z9").Value = 0
$ws.Range("AA9").Value = 0
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = 0
$ws.Range("AD9").Value = 0
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 0
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").Value = 0
$ws.Range("AJ9").Value = 0
$ws.Range("AK9").Value = 0
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = 0
$ws.Range("AN9").Value = 0
$ws.Range("AO9").Value = 0

# Row 8 now becomes a different fixture: Truro City vs Brackley Town.
$ws.Range("D8").Value = "Truro City"
$ws.Range("E8").Value = "Brackley Town"
$ws.Range("F8").Value = 2.2
$ws.Range("G8").Value = 2.98
$ws.Range("H8").Value = 2.66
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.4
$ws.Range("K8").Value = 5.9
$ws.Range("P8").Value = 1.7
$ws.Range("Q8").Value = 1.86
